$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (TestCaseName) + Column D (URI) edited first, in the exact
# order that reproduces the original authoring tool's shared-string slot
# reuse: D2 grabs the freed "TestMethod1" slot, A2 grabs the freed
# "TestMethod2" slot, A3 grabs the freed "http://facebook.com/" slot.
$ws.Range("D2").Value = "http://jsonplaceholder.typicode.com/"
$ws.Range("A2").Value = "validatingCommentsOfAParticularID"
$ws.Range("A3").Value = "ValidatingCommentWithPostID"

# New rows 4-11, column A test case names (appended as new shared strings).
$ws.Range("A4").Value = "responseStatusCodeVerificationForTodos"
$ws.Range("A5").Value = "responseStatusCodeVerificationForPhotos"
$ws.Range("A6").Value = "responseStatusCodeVerificationForAlbums"
$ws.Range("A7").Value = "responseStatusCodeVerificationForPosts"
$ws.Range("A8").Value = "responseStatusCodeVerificationForComments"
$ws.Range("A9").Value = "urlValidationOfaPhoto"
$ws.Range("A10").Value = "urlValidationOfaPhotoParameterized"
$ws.Range("A11").Value = "verifyPostedResource"

# Column B (Execute) and Column C (Environment) for every data row.
$ws.Range("B3").Value = "YES"
$ws.Range("C3").Value = "TEST"
for ($r = 4; $r -le 11; $r++) {
    $ws.Range("B" + $r).Value = "YES"
    $ws.Range("C" + $r).Value = "TEST"
}

# Column D (URI) for rows 3-11 — same shared URL string as D2.
for ($r = 3; $r -le 11; $r++) {
    $ws.Range("D" + $r).Value = "http://jsonplaceholder.typicode.com/"
}

# Hyperlinks for D2:D11 (all pointing at the same URL).
$ws.Hyperlinks.Add($ws.Range("D2"), "http://jsonplaceholder.typicode.com/")
$ws.Hyperlinks.Add($ws.Range("D3"), "http://jsonplaceholder.typicode.com/")
$ws.Hyperlinks.Add($ws.Range("D4"), "http://jsonplaceholder.typicode.com/")
$ws.Hyperlinks.Add($ws.Range("D6"), "http://jsonplaceholder.typicode.com/")
$ws.Hyperlinks.Add($ws.Range("D8"), "http://jsonplaceholder.typicode.com/")
$ws.Hyperlinks.Add($ws.Range("D5"), "http://jsonplaceholder.typicode.com/")
$ws.Hyperlinks.Add($ws.Range("D7"), "http://jsonplaceholder.typicode.com/")
$ws.Hyperlinks.Add($ws.Range("D9"), "http://jsonplaceholder.typicode.com/")
$ws.Hyperlinks.Add($ws.Range("D10"), "http://jsonplaceholder.typicode.com/")
$ws.Hyperlinks.Add($ws.Range("D11"), "http://jsonplaceholder.typicode.com/")

# Re-apply the Hyperlink cell style (Hyperlinks.Add swaps in a brand-new
# style slot; putting the built-in "Hyperlink" style back keeps every
# linked cell on the same style index the workbook already had).
$ws.Range("D2:D11").Style = "Hyperlink"

# Column A width grows to fit the long new test-case names.
$ws.Columns.Item(1).AutoFit()

# Selection moves to A3, matching the saved view state.
$ws.Range("A3").Select()
